$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update product name (B1) and short name (B2) on the ProductLoanInput sheet
$wsInput.Range("B1").Value = "4265-MS-EI-DB-DL-REC-RNI-FEE+INT-FFConMONTHLYonDAY25-FIFC-1-FFROP-DL-FIFR-1-MD-1st"
$wsInput.Range("B2").Value = "426v"

# Make ProductLoanInput the active sheet/tab (it previously pointed at
# ProductLoanOutput) and move the selection to B6, removing the test's
# dependency on another sheet being active.
$wsInput.Activate()
$wsInput.Range("B6").Select()
